# Apply crypto price/volume updates for Sat Jun 22 03:50:12 UTC 2024 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.277.54'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '3.505.15'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '583.63'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.07'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.91%  '
$ws.Range('D7').Value = '3.505.16'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.487'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E10').Value = '  +0.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.11'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  -3.46%  '
$ws.Range('D13').Value = '4.102.35'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000180'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D16').Value = '3.507.08'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.36'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -5.06%  '
$ws.Range('D18').Value = '64.298.69'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.74'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.82'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.21%  '
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '383.69'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.89%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.568'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.58%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.643.97'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '73.96'
$ws.Range('D25').Style = "Normal"
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('E28').Value = '  +4.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.59'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.57'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.28'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.86%  '
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('D34').Value = '3.522.63'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '23.55'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.35'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.18%  '
$ws.Range('E39').Value = '  -3.45%  '
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '164.23'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.47%  '
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '25.76'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '41.83'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.41'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('D50').Value = '2.473.39'
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.922'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.70%  '
